$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = 44268.50347222222
$ws.Range("J2").Value = 44275.50347222222
$ws.Range("K2").Value = 7
$ws.Range("I3").Value = 44284.45138888889
$ws.Range("J3").Value = 44291.45138888889
$ws.Range("K3").Value = 7
$ws.Range("I4").Value = 44296.19444444445
$ws.Range("J4").Value = 44303.19444444445
$ws.Range("K4").Value = 7
$ws.Range("I5").Value = 44324.625
$ws.Range("J5").Value = 44331.625
$ws.Range("K5").Value = 7
$ws.Range("I6").Value = 44335.53819444445
$ws.Range("J6").Value = 44335.65625
$ws.Range("K6").Value = 0.1180555555555556
$ws.Range("D7").Value = 24853
$ws.Range("F7").Value = 2744.51
$ws.Range("G7").Value = -94.44900000000052
$ws.Range("H7").Value = -0.005768719479134932
$ws.Range("I7").Value = 44335.65625
$ws.Range("J7").Value = 44342.65625
$ws.Range("K7").Value = 7
$ws.Range("B8").Value = 7
$ws.Range("G8").Value = -3102.68882
$ws.Range("I8").Value = 44367.89583333334
$ws.Range("J8").Value = 44374.89583333334
$ws.Range("K8").Value = 7
$ws.Range("B9").Value = 7
$ws.Range("G9").Value = 2636.770010000001
$ws.Range("I9").Value = 44374.97916666666
$ws.Range("J9").Value = 44381.97916666666
$ws.Range("K9").Value = 7
$ws.Range("B10").Value = 9
$ws.Range("G10").Value = 3918.140550000005
$ws.Range("I10").Value = 44398.22569444445
$ws.Range("J10").Value = 44405.22569444445
$ws.Range("K10").Value = 7
$ws.Range("B11").Value = 8
$ws.Range("G11").Value = 5011.852480000001
$ws.Range("I11").Value = 44412.59375
$ws.Range("J11").Value = 44419.59375
$ws.Range("K11").Value = 7
$ws.Range("I12").Value = 44440.72569444445
$ws.Range("J12").Value = 44447.72569444445
$ws.Range("K12").Value = 7
$ws.Range("B13").Value = 7
$ws.Range("G13").Value = -589.6995999999999
$ws.Range("I13").Value = 44451.39930555555
$ws.Range("J13").Value = 44453.63194444445
$ws.Range("K13").Value = 2.232638888888889
$ws.Range("B14").Value = 7
$ws.Range("G14").Value = 3096.922990000004
$ws.Range("I14").Value = 44470.44097222222
$ws.Range("J14").Value = 44477.44097222222
$ws.Range("K14").Value = 7
$ws.Range("B15").Value = 6
$ws.Range("G15").Value = 423.0571200000004
$ws.Range("I15").Value = 44489.56944444445
$ws.Range("J15").Value = 44496.56944444445
$ws.Range("K15").Value = 7
$ws.Range("B16").Value = 6
$ws.Range("G16").Value = 2508.688200000001
$ws.Range("I16").Value = 44526.36111111111
$ws.Range("J16").Value = 44533.36111111111
$ws.Range("K16").Value = 7
$ws.Range("B17").Value = 8
$ws.Range("G17").Value = 2906.980960000004
$ws.Range("I17").Value = 44534.22916666666
$ws.Range("J17").Value = 44541.22916666666
$ws.Range("K17").Value = 7
$ws.Range("B18").Value = -8
$ws.Range("D18").Value = 85307
$ws.Range("F18").Value = 3997.49
$ws.Range("G18").Value = -76.57528000000093
$ws.Range("H18").Value = -0.002400227332653104
$ws.Range("I18").Value = 44545.83680555555
$ws.Range("J18").Value = 44552.83680555555
$ws.Range("K18").Value = 7
$ws.Range("B19").Value = 10
$ws.Range("G19").Value = 540.9149000000025
$ws.Range("I19").Value = 44568.17708333334
$ws.Range("J19").Value = 44575.17708333334
$ws.Range("K19").Value = 7
$ws.Range("B20").Value = 14
$ws.Range("G20").Value = 3782.047640000006
$ws.Range("I20").Value = 44585.43402777778
$ws.Range("J20").Value = 44592.43402777778
$ws.Range("K20").Value = 7
$ws.Range("B21").Value = 13
$ws.Range("G21").Value = 3656.840460000005
$ws.Range("I21").Value = 44596.32291666666
$ws.Range("J21").Value = 44603.32291666666
$ws.Range("K21").Value = 7
$ws.Range("B22").Value = 16
$ws.Range("G22").Value = 7970.080000000009
$ws.Range("I22").Value = 44616.15625
$ws.Range("J22").Value = 44623.15625
$ws.Range("K22").Value = 7
$ws.Range("B23").Value = 15
$ws.Range("G23").Value = 5226.644850000001
$ws.Range("I23").Value = 44642.1875
$ws.Range("J23").Value = 44649.1875
$ws.Range("K23").Value = 7

$ws.Range("I2:J23").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("K2:K23").NumberFormat = "0"
